$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the invoice/list date (A1) by one month: 2024-04-24 -> 2024-05-24
$ws.Range("A1").Value = 45436

# Update the prices in column D for rows 28-31
$ws.Range("D28").Value = 230.1
$ws.Range("D29").Value = 300
$ws.Range("D30").Value = 336
$ws.Range("D31").Value = 422
